$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RTM")

# --- Column F ("Requirement model") header: mark the header cell with the
#     same "text" number format used for the new data cells below it. ---
$ws.Range("F1").NumberFormat = "@"

# --- New URS-version cells in column F that must be stored as literal TEXT
#     (their trailing zero would otherwise be lost if stored as a number,
#     e.g. "2.20" -> 2.2). Written in ascending value order so the
#     workbook's shared-string table is populated 2.20, 2.21, ..., 2.26 in
#     that order, matching the target layout. ---
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = "2.20"

$ws.Range("F23").NumberFormat = "@"
$ws.Range("F23").Value = "2.21"

$ws.Range("F20").NumberFormat = "@"
$ws.Range("F20").Value = "2.22"

$ws.Range("F18").NumberFormat = "@"
$ws.Range("F18").Value = "2.23"

$ws.Range("F19").NumberFormat = "@"
$ws.Range("F19").Value = "2.24"

$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "2.25"

$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "2.26"

# --- New URS-version cells in column F that are plain numbers (no trailing
#     zero to preserve), formatted the same as the rest of the column. ---
$ws.Range("F2").Value = 2.2
$ws.Range("F2").NumberFormat = "@"

$ws.Range("F3").Value = 2.2
$ws.Range("F3").NumberFormat = "@"

$ws.Range("F4").Value = 2.6
$ws.Range("F4").NumberFormat = "@"

$ws.Range("F5").Value = 2.14
$ws.Range("F5").NumberFormat = "@"

$ws.Range("F6").Value = 2.3
$ws.Range("F6").NumberFormat = "@"

$ws.Range("F7").Value = 2.4
$ws.Range("F7").NumberFormat = "@"

$ws.Range("F8").Value = 2.13
$ws.Range("F8").NumberFormat = "@"

$ws.Range("F9").Value = 2.6
$ws.Range("F9").NumberFormat = "@"

$ws.Range("F10").Value = 2.1
$ws.Range("F10").NumberFormat = "@"

$ws.Range("F11").Value = 2.7
$ws.Range("F11").NumberFormat = "@"

$ws.Range("F12").Value = 2.15
$ws.Range("F12").NumberFormat = "@"

$ws.Range("F13").Value = 2.16
$ws.Range("F13").NumberFormat = "@"

$ws.Range("F16").Value = 2.17
$ws.Range("F16").NumberFormat = "@"

$ws.Range("F17").Value = 2.18
$ws.Range("F17").NumberFormat = "@"

$ws.Range("F21").Value = 2.19
$ws.Range("F21").NumberFormat = "@"

# --- Row 9 (use case #8): reassigned from "Maintain station" to
#     "Maintain user", with its linked test-case count updated 6 -> 4. ---
$ws.Range("E9").Value = "Maintain user"
$ws.Range("G9").Value = 4
